$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.979.90"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "2.169.18"
$ws.Range("E3").Value = "  -2.39%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'246.51"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("E6").Value = "  -1.82%  "
$ws.Range("D7").Value = "'66.34"
$ws.Range("E7").Value = "  -6.07%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").Value = "'58.45"
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").Value = "'0.0926"
$ws.Range("E11").Value = "  -3.45%  "
$ws.Range("D12").Value = "'35.65"
$ws.Range("E12").Value = "  -15.33%  "
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("D14").Value = "'6.89"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").Value = "2.490.59"
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("D16").Value = "'0.858"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").Value = "'14.28"
$ws.Range("E17").Value = "  -4.05%  "
$ws.Range("D18").Value = "2.162.71"
$ws.Range("E18").Value = "  -2.50%  "
$ws.Range("D19").Value = "40.886.33"
$ws.Range("E19").Value = "  -1.36%  "
$ws.Range("D20").Value = "0.0₃0939"
$ws.Range("E20").Value = "  -2.56%  "
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("D22").Value = "'71.29"
$ws.Range("E22").Value = "  -2.13%  "
$ws.Range("D23").Value = "'229.10"
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("E24").Value = "  -6.03%  "
$ws.Range("D25").Value = "'11.47"
$ws.Range("E25").Value = "  +12.23%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'3.71"
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("E28").Value = "  -3.63%  "
$ws.Range("E29").Value = "  -5.64%  "
$ws.Range("D30").Value = "'168.92"
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("E31").Value = "  -8.58%  "
$ws.Range("D32").Value = "'20.17"
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("D33").Value = "'0.121"
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("D34").Value = "'5.71"
$ws.Range("E34").Value = "  +3.48%  "
$ws.Range("D35").Value = "'0.0747"
$ws.Range("E35").Value = "  +4.35%  "
$ws.Range("E36").Value = "  -2.77%  "
$ws.Range("E37").Value = "  -2.54%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "'24.96"
$ws.Range("E38").Value = "  -5.01%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'3.99"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("D40").Value = "'0.0299"
$ws.Range("E40").Value = "  +5.09%  "
$ws.Range("D41").Value = "'2.18"
$ws.Range("E41").Value = "  -4.93%  "
$ws.Range("D42").Value = "'5.46"
$ws.Range("E42").Value = "  -8.89%  "
$ws.Range("E43").Value = "  -3.71%  "
$ws.Range("D44").Value = "'60.31"
$ws.Range("E44").Value = "  -12.64%  "
$ws.Range("D45").Value = "'4.82"
$ws.Range("E45").Value = "  -4.98%  "
$ws.Range("E46").Value = "  -7.65%  "
$ws.Range("D47").Value = "'8.45"
$ws.Range("E47").Value = "  -3.27%  "
$ws.Range("D48").Value = "'0.0991"
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("E51").Value = "  -3.15%  "
